$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Header 2 holds the "BTec_Logo-Orange" picture, currently exported as
# image1.jpg -- rename it to image2.jpg.
$h2 = $sec.Headers.Item(2)
if ($h2.Exists) {
    for ($si = 1; $si -le $h2.Range.InlineShapes.Count; $si++) {
        $s = $h2.Range.InlineShapes.Item($si)
        if ($s.AlternativeText -eq "BTec_Logo-Orange") {
            $s.Name = "image2.jpg"
        }
    }
}

# Footer 1 and Footer 2 both hold the Pearson Edexcel logo picture,
# currently exported as image2.png -- rename it to image1.png.
for ($fi = 1; $fi -le $sec.Footers.Count; $fi++) {
    $f = $sec.Footers.Item($fi)
    if ($f.Exists) {
        for ($si = 1; $si -le $f.Range.InlineShapes.Count; $si++) {
            $s = $f.Range.InlineShapes.Item($si)
            if ($s.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $s.Name = "image1.png"
            }
        }
    }
}
